$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray <w:bookmarkStart/.../w:bookmarkEnd> for "_GoBack"
#    that currently sits at the top of the SMA paragraph (paragraph 3).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Turn the trailing empty paragraph into a new "Sharpe ratio: <link>"
#    paragraph, followed by a (new) empty paragraph - i.e. insert the
#    Sharpe-ratio paragraph just before the final empty paragraph.
# ---------------------------------------------------------------------
$last = $d.Paragraphs.Last
$insertRange = $last.Range
$insertRange.Collapse(1)  # wdCollapseStart - insert before the final empty paragraph

$sharpeXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Sharpe</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>ratio</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve">: </w:t></w:r>
<w:hyperlink r:id="rIdSharpeRatio" w:history="1">
<w:r><w:t>http://www.investopedia.com/terms/s/sharperatio.asp</w:t></w:r>
</w:hyperlink>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
<Relationship Id="rIdSharpeRatio" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="http://www.investopedia.com/terms/s/sharperatio.asp" TargetMode="External"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertRange.InsertXML($sharpeXml)

# The hyperlink display text needs the character style "Hyperlink"
# (InsertXML drops w:rStyle references, so (re)apply it explicitly).
$hyperlinkCount = $d.Hyperlinks.Count
$newHyperlink = $d.Hyperlinks.Item($hyperlinkCount)
$hlStart = $newHyperlink.Range.Start
$hlEnd = $newHyperlink.Range.End
$hlRange = $d.Range($hlStart, $hlEnd)
$hlRange.Style = "Hyperlink"

# ---------------------------------------------------------------------
# 3) Put the "_GoBack" bookmark back on the (now) final empty paragraph.
# ---------------------------------------------------------------------
$finalPara = $d.Paragraphs.Last
$d.Bookmarks.Add("_GoBack", $finalPara.Range)

Write-Output $d.Content.Text
